# Apply the changes described in the commit:
#  - employee quarterly rates (value corrections to existing short names)
#  - names of functional groups (new highlighted reference rows)
#  - hourly production calendar / other new reference rows
#  - highlight the new "functional group" rows with an orange fill

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First, grow the table by one row, copying the formatting of the last
#     (still blank) placeholder row down to the new blank row 39 ---
$ws.Range("A38:B38").Copy()
$ws.Range("A39:B39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fix two existing short-name values (quarterly-rate corrections) ---
$ws.Range("B6").Value = "ДПИС"

$ws.Range("B9").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value = "MES"

# --- New plain rows (25-27), formatted like the existing data rows ---
$ws.Range("A2").Copy()
$ws.Range("A25:A27").PasteSpecial(-4122)
$ws.Range("B9").Copy()
$ws.Range("B25:B27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A25").Value = "ООО ""АрСи БиАй"""
$ws.Range("B25").Value = "MES"
$ws.Range("A26").Value = "Системы управления производством"
$ws.Range("B26").Value = "MES"
$ws.Range("A27").Value = "Отдел сопровождения систем геологоразведки и добычи"
$ws.Range("B27").Value = "ГиД"

# --- New rows (28-37): functional-group names, column A highlighted orange ---
$ws.Range("A28").Interior.Color = 49407
$ws.Range("A28").Copy()
$ws.Range("A29:A37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B9").Copy()
$ws.Range("B28:B30").PasteSpecial(-4122)
$ws.Range("B32:B33").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("B31").PasteSpecial(-4122)
$ws.Range("B34:B37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A28").Value = "Системы контроля качества (СКК)"
$ws.Range("B28").Value = "СКК"
$ws.Range("A29").Value = "Системы управления производством (MES)"
$ws.Range("B29").Value = "MES"
$ws.Range("A30").Value = "Системы промышленной безопасности и экологии (HSE)"
$ws.Range("B30").Value = "СПБиЭ"
$ws.Range("A31").Value = "Производственные системы геологоразведки и добычи (ГиД)"
$ws.Range("B31").Value = "ГиД"
$ws.Range("A32").Value = "Корпоративные информационные системы (КИС)"
$ws.Range("B32").Value = "КИС"
$ws.Range("A33").Value = "Портальные решения (ПР)"
$ws.Range("B33").Value = "Порталы"
$ws.Range("A34").Value = "Системы управления персоналом (НСМ)"
$ws.Range("B34").Value = "HCM"
$ws.Range("A35").Value = "Системы управления предприятием (ERP)"
$ws.Range("B35").Value = "ERP"
$ws.Range("A36").Value = "Системы бизнес-анализа (BI)"
$ws.Range("B36").Value = "BI"
$ws.Range("A37").Value = "Профессиональные услуги (ПУ)"
$ws.Range("B37").Value = "ПУ"

# --- New plain row (38) ---
$ws.Range("A2").Copy()
$ws.Range("A38").PasteSpecial(-4122)
$ws.Range("B38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A38").Value = "Департамент поддержки информационных систем"
$ws.Range("B38").Value = "ДПИС"

# --- Move selection to reflect where the editor last worked ---
[void]$ws.Range("B26").Select()
